$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.153101325035095
$ws.Range("B1").Value = 2.258035659790039
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.140463352203369
$ws.Range("E1").Value = 1.060437679290771
